$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the AST column (D) as tested ("y") for the shift instructions:
# LSL (row 16), LSR (row 17), ASR (row 18), ROR (row 19), RRX (row 20)
$ws.Range("D16").Value = "y"
$ws.Range("D17").Value = "y"
$ws.Range("D18").Value = "y"
$ws.Range("D19").Value = "y"
$ws.Range("D20").Value = "y"

# Update the active selection to D20, matching where the edit was made
$ws.Range("D20").Select()
